$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.Value = "'" + $value
    $cell.Style = $origStyle
}

# Row 2
$ws.Range("D2").Value = "29.301.32"
$ws.Range("E2").Value = "  -0.08%  "

# Row 3
$ws.Range("D3").Value = "1.832.03"
$ws.Range("E3").Value = "  -0.43%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.005"
$ws.Range("E4").Value = "  +0.44%  "

# Row 5
Set-TextValue $ws.Range("D5") "235.85"
$ws.Range("E5").Value = "  -1.48%  "

# Row 6
Set-TextValue $ws.Range("D6") "0.6038"
$ws.Range("E6").Value = "  -2.90%  "

# Row 7
$ws.Range("E7").Value = "  +0.32%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.06986"
$ws.Range("E8").Value = "  -4.90%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.2776"
$ws.Range("E9").Value = "  -3.47%  "

# Row 10
Set-TextValue $ws.Range("D10") "23.63"
$ws.Range("E10").Value = "  -4.10%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.07616"
$ws.Range("E11").Value = "  -1.26%  "

# Row 12
$ws.Range("D12").Value = "1.848.68"
$ws.Range("E12").Value = "  +0.74%  "

# Row 13
Set-TextValue $ws.Range("D13") "4.768"
$ws.Range("E13").Value = "  -3.45%  "

# Row 14
Set-TextValue $ws.Range("D14") "0.6313"
$ws.Range("E14").Value = "  -4.20%  "

# Row 15
Set-TextValue $ws.Range("D15") "0.000009892"
$ws.Range("E15").Value = "  -6.09%  "

# Row 16
Set-TextValue $ws.Range("D16") "77.95"
$ws.Range("E16").Value = "  -4.14%  "

# Row 17
$ws.Range("D17").Value = "28.951.52"
$ws.Range("E17").Value = "  -1.25%  "

# Row 18
Set-TextValue $ws.Range("D18") "5.600"
$ws.Range("E18").Value = "  -10.06%  "

# Row 19
Set-TextValue $ws.Range("D19") "218.25"
$ws.Range("E19").Value = "  -7.64%  "

# Row 20
Set-TextValue $ws.Range("D20") "1.005"
$ws.Range("E20").Value = "  +0.44%  "

# Row 21
Set-TextValue $ws.Range("D21") "11.60"
$ws.Range("E21").Value = "  -4.72%  "

# Row 22
Set-TextValue $ws.Range("D22") "6.918"
$ws.Range("E22").Value = "  -3.79%  "

# Row 23
Set-TextValue $ws.Range("D23") "1.003"
$ws.Range("E23").Value = "  -0.35%  "

# Row 24
Set-TextValue $ws.Range("D24") "156.33"
$ws.Range("E24").Value = "  -0.64%  "

# Row 25
Set-TextValue $ws.Range("D25") "7.994"
$ws.Range("E25").Value = "  -4.80%  "

# Row 26
Set-TextValue $ws.Range("D26") "0.1296"
$ws.Range("E26").Value = "  -2.62%  "

# Row 27
Set-TextValue $ws.Range("D27") "16.56"
$ws.Range("E27").Value = "  -3.74%  "

# Row 28
Set-TextValue $ws.Range("D28") "0.06487"
$ws.Range("E28").Value = "  -5.36%  "

# Row 29
Set-TextValue $ws.Range("D29") "1.427"
$ws.Range("E29").Value = "  -3.19%  "

# Row 30
Set-TextValue $ws.Range("D30") "1.445"
$ws.Range("E30").Value = "  -2.25%  "

# Row 31
Set-TextValue $ws.Range("D31") "3.852"
$ws.Range("E31").Value = "  -1.87%  "

# Row 32
Set-TextValue $ws.Range("D32") "3.801"
$ws.Range("E32").Value = "  -5.14%  "

# Row 33
Set-TextValue $ws.Range("D33") "1.738"
$ws.Range("E33").Value = "  -0.26%  "

# Row 34
Set-TextValue $ws.Range("D34") "1.099"
$ws.Range("E34").Value = "  -4.52%  "

# Row 35
Set-TextValue $ws.Range("D35") "0.6509"
$ws.Range("E35").Value = "  -4.11%  "

# Row 36
Set-TextValue $ws.Range("D36") "2.543"
$ws.Range("E36").Value = "  -1.50%  "

# Row 37
Set-TextValue $ws.Range("D37") "2.763"
$ws.Range("E37").Value = "  -0.64%  "

# Row 38
$ws.Range("E38").Value = "  -3.20%  "

# Row 39
Set-TextValue $ws.Range("D39") "6.574"
$ws.Range("E39").Value = "  -0.94%  "

# Row 40
$ws.Range("D40").Value = "1.148.44"
$ws.Range("E40").Value = "  -6.72%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.8954"
$ws.Range("E41").Value = "  -5.00%  "

# Row 42
Set-TextValue $ws.Range("D42") "1.004"
$ws.Range("E42").Value = "  +0.23%  "

# Row 43
$ws.Range("D43").Value = "1.995.64"
$ws.Range("E43").Value = "  +0.36%  "

# Row 44
Set-TextValue $ws.Range("D44") "100.95"
$ws.Range("E44").Value = "  -0.30%  "

# Row 45
Set-TextValue $ws.Range("D45") "62.37"
$ws.Range("E45").Value = "  -3.94%  "

# Row 46
Set-TextValue $ws.Range("D46") "0.00000000114"
$ws.Range("E46").Value = "  -4.39%  "

# Row 47
Set-TextValue $ws.Range("D47") "1.624"
$ws.Range("E47").Value = "  -3.50%  "

# Row 48
Set-TextValue $ws.Range("D48") "8.566"
$ws.Range("E48").Value = "  -2.78%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.4552"
$ws.Range("E49").Value = "  -0.37%  "

# Row 50
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D50") "6.441"
$ws.Range("E50").Value = "  -6.11%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D51") "0.05497"
$ws.Range("E51").Value = "  -2.33%  "
